$wb = $excel.ActiveWorkbook

$wsIndisp = $wb.Worksheets.Item("Persone-Indisp")
$wsTurni  = $wb.Worksheets.Item("Turni Fissi")

# --- Update the "Turni Fissi" (fixed shifts) assignments ---
# Night-shift detail cells (row 3-5, columns F/D/C/G)
$wsTurni.Range("F3").Value = "MAR"
$wsTurni.Range("D4").Value = "CMG"
$wsTurni.Range("C5").Value = "MAR"
$wsTurni.Range("G5").Value = "CMG"

# "Ricerca" column (B) re-assignments across several day-blocks
$wsTurni.Range("B6").Value = "CAR"

$wsTurni.Range("B13").Value = "MAR"
$wsTurni.Range("B14").Value = "MAR"
$wsTurni.Range("B15").Value = "MAR"
$wsTurni.Range("B16").Value = "MAR"
$wsTurni.Range("B17").Value = "MAR"

$wsTurni.Range("B20").Value = "MAD"
$wsTurni.Range("B21").Value = "MAD"
$wsTurni.Range("B22").Value = "MAD"
$wsTurni.Range("B23").Value = "MAD"
$wsTurni.Range("B24").Value = "MAD"

$wsTurni.Range("B27").Value = "DAN"
$wsTurni.Range("B28").Value = "DAN"
$wsTurni.Range("B29").Value = "DAN"
$wsTurni.Range("B30").Value = "DAN"
$wsTurni.Range("B31").Value = "DAN"

# --- Final view state: user ends on "Persone-Indisp" selecting AC6, then
#     switches to and finishes on "Turni Fissi" selecting J3 ---
$wsIndisp.Activate()
$wsIndisp.Range("AC6").Select()

$wsTurni.Activate()
$wsTurni.Range("J3").Select()
